$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns D (VAR(OK)), E (MSPE), F (S_nugget), H (VAR(DATA))
# for rows 2 through 10, reflecting the optimized variogram parameters for area4.

$data = @(
    @{ Row = 2;  D = 0.6499262231420104; E = 6.032319885138804;  F = 0.8921;   H = 6.657481653811705 },
    @{ Row = 3;  D = 0.9074517110996223; E = 5.911469804138336;  F = 1.8169;   H = 6.657481653811705 },
    @{ Row = 4;  D = 1.127206373884514;  E = 5.807892747486808;  F = 2.2282;   H = 6.657481653811705 },
    @{ Row = 5;  D = 1.267178822006443;  E = 5.69176688981334;   F = 2.6597;   H = 6.657481653811705 },
    @{ Row = 6;  D = 1.375422044123241;  E = 5.606584187100444;  F = 2.9632;   H = 6.657481653811705 },
    @{ Row = 7;  D = 1.447466125483911;  E = 5.582565396729682;  F = 3.2233;   H = 6.657481653811705 },
    @{ Row = 8;  D = 1.483948129415954;  E = 5.564613845110005;  F = 3.4306;   H = 6.657481653811705 },
    @{ Row = 9;  D = 1.506771102985734;  E = 5.556891163051717;  F = 3.5275;   H = 6.657481653811705 },
    @{ Row = 10; D = 1.53187157506313;   E = 5.555116743832956;  F = 3.60112;  H = 6.657481653811705 }
)

foreach ($row in $data) {
    $ws.Range("D$($row.Row)").Value = $row.D
    $ws.Range("E$($row.Row)").Value = $row.E
    $ws.Range("F$($row.Row)").Value = $row.F
    $ws.Range("H$($row.Row)").Value = $row.H
}
